$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.657.53'
$ws.Range("E2").Value = '  -0.46%  '

$ws.Range("D3").Value = '1.597.46'
$ws.Range("E3").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.34'
$ws.Range("E5").Value = '  +0.35%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("E8").Value = '  +0.10%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.245'
$ws.Range("E9").Value = '  -1.08%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.69'
$ws.Range("E10").Value = '  +0.32%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0837'
$ws.Range("E11").Value = '  -0.02%  '

$ws.Range("D12").Value = '1.820.69'
$ws.Range("E12").Value = '  -0.04%  '

$ws.Range("D13").Value = '1.633.04'
$ws.Range("E13").Value = '  +1.21%  '

$ws.Range("E14").Value = '  -0.41%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.522'
$ws.Range("E15").Value = '  -1.34%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.87'
$ws.Range("E16").Value = '  +2.26%  '

$ws.Range("D17").Value = '26.653.47'
$ws.Range("E17").Value = '  -0.33%  '

$ws.Range("D18").Value = '0.0₃0728'
$ws.Range("E18").Value = '  -0.07%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '209.66'
$ws.Range("E19").Value = '  +0.37%  '

$ws.Range("E20").Value = '  +0.02%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.78'
$ws.Range("E21").Value = '  +0.74%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.27'
$ws.Range("E22").Value = '  -0.14%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.29'
$ws.Range("E23").Value = '  -1.43%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.92'
$ws.Range("E24").Value = '  +0.85%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.25'
$ws.Range("E25").Value = '  -0.08%  '

$ws.Range("E26").Value = '  -0.10%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.19'
$ws.Range("E27").Value = '  -4.16%  '

$ws.Range("E28").Value = '  +2.53%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.29'
$ws.Range("E29").Value = '  +0.07%  '

$ws.Range("E30").Value = '  +0.78%  '

$ws.Range("E31").Value = '  +0.09%  '

$ws.Range("E32").Value = '  -0.63%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.666'
$ws.Range("E33").Value = '  -0.40%  '

$ws.Range("E34").Value = '  -0.76%  '

$ws.Range("D35").Value = '1.301.18'
$ws.Range("E35").Value = '  -0.92%  '

$ws.Range("E36").Value = '  +0.70%  '

$ws.Range("E37").Value = '  -2.15%  '

$ws.Range("E38").Value = '  -0.95%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.841'
$ws.Range("E39").Value = '  +2.61%  '

$ws.Range("E40").Value = '  +0.03%  '

$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.38'
$ws.Range("E41").Value = '  +2.18%  '

$ws.Range("B42").Value = 'MXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.20'
$ws.Range("E42").Value = '  +1.23%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.788'
$ws.Range("E43").Value = '  +0.16%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '63.82'
$ws.Range("E44").Value = '  +1.56%  '

$ws.Range("D45").Value = '1.733.64'
$ws.Range("E45").Value = '  -0.08%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.893'
$ws.Range("E46").Value = '  +10.06%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.22'
$ws.Range("E47").Value = '  +1.39%  '

$ws.Range("E48").Value = '  +0.92%  '

$ws.Range("E49").Value = '  +2.81%  '

$ws.Range("E50").Value = '  -0.91%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.47'
$ws.Range("E51").Value = '  +0.60%  '
